$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Duel Decks Merfolk vs. Goblins Tokens (TDDT)"
$ws.Range("A2").Value = "Elemental Shaman"

$ws.Range("A3").Delete()
$ws.Range("A3").Delete()
